$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: update title/link
$ws.Range("D9").Value = "온라인 교육 시스템이 미래인 이유"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/online-education-the-future/#utm_source=rss&utm_medium=rss&utm_campaign=online-education-the-future"

# Row 46: update title/link
$ws.Range("D46").Value = "[Bioinformatics] 2021년 12월,  한국유전체학회 제18회 동계 워크샵"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/422"
